# Update "想去人数" (F) and "最低票价" (G) figures on the two sheets that
# carry the full event list: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 78
    $ws.Range("F3").Value = 69
    $ws.Range("F4").Value = 3734
    $ws.Range("F5").Value = 2262
    $ws.Range("F6").Value = 442
    $ws.Range("F7").Value = 7
    $ws.Range("F8").Value = 13

    if ($sheetName -eq "展览") {
        $ws.Range("F11").Value = 83
        $ws.Range("F12").Value = 1381
        $ws.Range("F14").Value = 2201
    } else {
        $ws.Range("F12").Value = 83
        $ws.Range("F15").Value = 1381
        $ws.Range("F17").Value = 2201
    }
}
